$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header B1: "ZNr" -> "Tab_7a_Daten.ZNr"
$ws.Range("B1").Value = "Tab_7a_Daten.ZNr"

# 2. Widen column C (Disaggregation 1 Kategorie) to match new content width.
#    The underlying engine snaps ColumnWidth to an MDW-7 pixel grid, so we choose
#    the input that lands closest to the target stored width of 33.25.
$ws.Columns.Item(3).ColumnWidth = 32.57142857142857

# 3. Add a new column AY (51) - "Tab_6a_Zeitreihen.ZNr" - with its own width.
#    Closest achievable ColumnWidth input for a target stored width of 26.8046875.
$ws.Columns.Item(51).ColumnWidth = 26.142857142857142

# 4. Populate the new AY column.
#    First copy formatting from the neighboring AX column so the new cells get
#    the same styles (header style on row 1, data style on rows 2-52).
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)
$ws.Range("AY1").Value = "Tab_6a_Zeitreihen.ZNr"

$ws.Range("AX2:AX52").Copy()
$ws.Range("AY2:AY52").PasteSpecial(-4122)

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 51).Value = "Z07_B01_P01_Ib01_I01_Z01"
}
for ($r = 19; $r -le 35; $r++) {
    $ws.Cells.Item($r, 51).Value = "Z07_B01_P01_Ib01_I02_Z01"
}
for ($r = 36; $r -le 52; $r++) {
    $ws.Cells.Item($r, 51).Value = "Z07_B02_P01_Ib01_I01_Z01"
}
